# Atualização de bases das ligas, do dia: 29-03-2024 às 13:24
# Updates a batch of match/odds rows (re-ordering of fixtures 130-133 and 142-145,
# plus odds refresh for rows 184-191) in the "Ecuador LigaPro Serie A" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range("B130").Value = 7483247
$ws.Range("F130").Value = "Mushuc Runa"
$ws.Range("G130").Value = "Universidad Catolica del Ecuador"
$ws.Range("I130").Value = 2
$ws.Range("J130").Value = "A"
$ws.Range("K130").Value = 3.25
$ws.Range("L130").Value = 3.2
$ws.Range("M130").Value = 2.25
$ws.Range("N130").Value = 3.5
$ws.Range("O130").Value = 3.25
$ws.Range("P130").Value = 2.1
$ws.Range("Q130").Value = 0.5
$ws.Range("R130").Value = 1.775
$ws.Range("S130").Value = 2.025
$ws.Range("U130").Value = 1.9
$ws.Range("V130").Value = 1.9
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 1.1
$ws.Range("AA130").Value = 1.025
$ws.Range("AC130").Value = 0.8999999999999999

# Row 131
$ws.Range("B131").Value = 7483081
$ws.Range("F131").Value = "Deportivo Cuenca"
$ws.Range("G131").Value = "El Nacional"
$ws.Range("H131").Value = 1
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 2.75
$ws.Range("L131").Value = 3.25
$ws.Range("M131").Value = 2.55
$ws.Range("N131").Value = 3
$ws.Range("O131").Value = 3.3
$ws.Range("P131").Value = 2.3
$ws.Range("Q131").Value = 0.25
$ws.Range("R131").Value = 1.825
$ws.Range("S131").Value = 1.975
$ws.Range("T131").Value = 2.75
$ws.Range("U131").Value = 2
$ws.Range("V131").Value = 1.8
$ws.Range("W131").Value = 2
$ws.Range("X131").Value = -1
$ws.Range("Z131").Value = 0.825
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = -1
$ws.Range("AC131").Value = 0.8

# Row 132
$ws.Range("B132").Value = 7483189
$ws.Range("F132").Value = "Independiente del Valle"
$ws.Range("G132").Value = "Orense"
$ws.Range("H132").Value = 2
$ws.Range("J132").Value = "D"
$ws.Range("K132").Value = 1.4
$ws.Range("L132").Value = 4.75
$ws.Range("M132").Value = 7
$ws.Range("N132").Value = 1.4
$ws.Range("O132").Value = 4.5
$ws.Range("P132").Value = 8
$ws.Range("Q132").Value = -1.25
$ws.Range("R132").Value = 1.875
$ws.Range("S132").Value = 1.925
$ws.Range("U132").Value = 1.925
$ws.Range("V132").Value = 1.875
$ws.Range("X132").Value = 3.5
$ws.Range("Y132").Value = -1
$ws.Range("AA132").Value = 0.925
$ws.Range("AB132").Value = 0.925
$ws.Range("AC132").Value = -1

# Row 133
$ws.Range("B133").Value = 7483281
$ws.Range("F133").Value = "SD Aucas"
$ws.Range("G133").Value = "Delfin SC"
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = "D"
$ws.Range("K133").Value = 1.909
$ws.Range("M133").Value = 4.2
$ws.Range("N133").Value = 1.909
$ws.Range("O133").Value = 3.5
$ws.Range("P133").Value = 4
$ws.Range("Q133").Value = -0.5
$ws.Range("R133").Value = 1.9
$ws.Range("S133").Value = 1.9
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.8
$ws.Range("V133").Value = 2
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 2.5
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.8999999999999999
$ws.Range("AC133").Value = 1

# Row 142
$ws.Range("B142").Value = 7528848
$ws.Range("F142").Value = "Emelec"
$ws.Range("G142").Value = "Deportivo Cuenca"
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = "H"
$ws.Range("K142").Value = 1.75
$ws.Range("L142").Value = 3.5
$ws.Range("M142").Value = 4.2
$ws.Range("N142").Value = 2.4
$ws.Range("O142").Value = 3.1
$ws.Range("P142").Value = 2.75
$ws.Range("R142").Value = 2.05
$ws.Range("S142").Value = 1.75
$ws.Range("U142").Value = 1.8
$ws.Range("V142").Value = 2
$ws.Range("W142").Value = 1.4
$ws.Range("X142").Value = -1
$ws.Range("Z142").Value = 1.05
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8

# Row 143
$ws.Range("B143").Value = 7528858
$ws.Range("F143").Value = "Orense"
$ws.Range("G143").Value = "SD Aucas"
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 2
$ws.Range("J143").Value = "A"
$ws.Range("K143").Value = 2.2
$ws.Range("L143").Value = 3.2
$ws.Range("M143").Value = 3.2
$ws.Range("N143").Value = 1.95
$ws.Range("O143").Value = 3.2
$ws.Range("P143").Value = 3.8
$ws.Range("Q143").Value = -0.5
$ws.Range("R143").Value = 1.95
$ws.Range("S143").Value = 1.85
$ws.Range("U143").Value = 1.85
$ws.Range("V143").Value = 1.95
$ws.Range("W143").Value = -1
$ws.Range("Y143").Value = 2.8
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.8500000000000001
$ws.Range("AB143").Value = 0.8500000000000001

# Row 144
$ws.Range("B144").Value = 7528857
$ws.Range("F144").Value = "Universidad Catolica del Ecuador"
$ws.Range("G144").Value = "Barcelona Guayaquil"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("K144").Value = 1.533
$ws.Range("L144").Value = 4
$ws.Range("M144").Value = 5.5
$ws.Range("N144").Value = 1.5
$ws.Range("O144").Value = 4.333
$ws.Range("P144").Value = 5.25
$ws.Range("Q144").Value = -1
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2
$ws.Range("T144").Value = 3
$ws.Range("U144").Value = 1.975
$ws.Range("V144").Value = 1.825
$ws.Range("Y144").Value = 4.25
$ws.Range("AA144").Value = 1
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.825

# Row 145
$ws.Range("B145").Value = 7528852
$ws.Range("F145").Value = "Delfin SC"
$ws.Range("G145").Value = "Tecnico Universitario"
$ws.Range("H145").Value = 2
$ws.Range("I145").Value = 2
$ws.Range("J145").Value = "D"
$ws.Range("K145").Value = 2.1
$ws.Range("L145").Value = 3.4
$ws.Range("M145").Value = 3.1
$ws.Range("N145").Value = 2.1
$ws.Range("O145").Value = 3.4
$ws.Range("P145").Value = 3.1
$ws.Range("Q145").Value = -0.25
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 1.9
$ws.Range("V145").Value = 1.9
$ws.Range("X145").Value = 2.4
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -0.5
$ws.Range("AA145").Value = 0.5
$ws.Range("AB145").Value = 0.8999999999999999
$ws.Range("AC145").Value = -1

# Row 184
$ws.Range("H184").Value = 5
$ws.Range("I184").Value = 0
$ws.Range("J184").Value = "H"
$ws.Range("N184").Value = 1.3
$ws.Range("O184").Value = 5.25
$ws.Range("P184").Value = 7.5
$ws.Range("R184").Value = 1.9
$ws.Range("S184").Value = 1.9
$ws.Range("U184").Value = 1.8
$ws.Range("V184").Value = 2
$ws.Range("W184").Value = 0.3
$ws.Range("X184").Value = -1
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = 0.8999999999999999
$ws.Range("AA184").Value = -1
$ws.Range("AB184").Value = 0.8
$ws.Range("AC184").Value = -1

# Row 185
$ws.Range("N185").Value = 2.3
$ws.Range("O185").Value = 3.25
$ws.Range("P185").Value = 2.75
$ws.Range("R185").Value = 1.75
$ws.Range("S185").Value = 2.05
$ws.Range("U185").Value = 2.025
$ws.Range("V185").Value = 1.775

# Row 186
$ws.Range("R186").Value = 1.925
$ws.Range("S186").Value = 1.875

# Row 187
$ws.Range("R187").Value = 1.975
$ws.Range("S187").Value = 1.825

# Row 188
$ws.Range("N188").Value = 4.333
$ws.Range("O188").Value = 3.6
$ws.Range("P188").Value = 1.727
$ws.Range("Q188").Value = 0.75
$ws.Range("R188").Value = 1.9
$ws.Range("S188").Value = 1.9
$ws.Range("U188").Value = 1.9
$ws.Range("V188").Value = 1.9

# Row 190
$ws.Range("N190").Value = 1.8
$ws.Range("O190").Value = 3.4
$ws.Range("P190").Value = 4.2
$ws.Range("Q190").Value = -0.5
$ws.Range("R190").Value = 1.8
$ws.Range("S190").Value = 2
$ws.Range("U190").Value = 1.9
$ws.Range("V190").Value = 1.9

# Row 191
$ws.Range("N191").Value = 2.3
$ws.Range("O191").Value = 3.25
$ws.Range("R191").Value = 1.975
$ws.Range("S191").Value = 1.825
$ws.Range("T191").Value = 2.5
$ws.Range("U191").Value = 1.975
$ws.Range("V191").Value = 1.825
